$wb = $excel.ActiveWorkbook

# ---- Sheet LP1912 (sheet1) ----
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = 'Última actualización: 20:31:05'
$ws1.Range("A3").Value = 'Total filas: 524'

$d1_119 = New-Object "object[,]" 2,5
$d1_119[0,0] = '08:57:13'
$d1_119[0,1] = '09:35'
$d1_119[0,2] = '23_HERNANDEZ'
$d1_119[0,3] = 38
$d1_119[0,4] = 'LP1912'
$d1_119[1,0] = '08:57:13'
$d1_119[1,1] = '09:35'
$d1_119[1,2] = '16_SANTA ANA'
$d1_119[1,3] = 38
$d1_119[1,4] = 'LP1912'
$ws1.Range("A119:E120").Value = $d1_119

$d1_201 = New-Object "object[,]" 3,5
$d1_201[0,0] = '10:26:41'
$d1_201[0,1] = '12:21'
$d1_201[0,2] = '26_HERNANDEZ'
$d1_201[0,3] = 115
$d1_201[0,4] = 'LP1912'
$d1_201[1,0] = '11:20:07'
$d1_201[1,1] = '12:21'
$d1_201[1,2] = '14_ABASTO'
$d1_201[1,3] = 61
$d1_201[1,4] = 'LP1912'
$d1_201[2,0] = '11:20:07'
$d1_201[2,1] = '12:21'
$d1_201[2,2] = '215A_EL PATO'
$d1_201[2,3] = 61
$d1_201[2,4] = 'LP1912'
$ws1.Range("A201:E203").Value = $d1_201

$d1_264 = New-Object "object[,]" 2,5
$d1_264[0,0] = '12:37:14'
$d1_264[0,1] = '14:20'
$d1_264[0,2] = '215C_EL PATO'
$d1_264[0,3] = 103
$d1_264[0,4] = 'LP1912'
$d1_264[1,0] = '13:19:56'
$d1_264[1,1] = '14:20'
$d1_264[1,2] = '26_HERNANDEZ'
$d1_264[1,3] = 61
$d1_264[1,4] = 'LP1912'
$ws1.Range("A264:E265").Value = $d1_264

$d1_299 = New-Object "object[,]" 2,5
$d1_299[0,0] = '14:19:48'
$d1_299[0,1] = '15:38'
$d1_299[0,2] = '215A_EL PATO'
$d1_299[0,3] = 79
$d1_299[0,4] = 'LP1912'
$d1_299[1,0] = '14:19:48'
$d1_299[1,1] = '15:38'
$d1_299[1,2] = '23_HERNANDEZ'
$d1_299[1,3] = 79
$d1_299[1,4] = 'LP1912'
$ws1.Range("A299:E300").Value = $d1_299

$d1_333 = New-Object "object[,]" 2,5
$d1_333[0,0] = '16:33:08'
$d1_333[0,1] = '16:34'
$d1_333[0,2] = '16_P MOR-SANTA ANA'
$d1_333[0,3] = 1
$d1_333[0,4] = 'LP1912'
$d1_333[1,0] = '15:57:48'
$d1_333[1,1] = '16:34'
$d1_333[1,2] = '23_HERNANDEZ'
$d1_333[1,3] = 37
$d1_333[1,4] = 'LP1912'
$ws1.Range("A333:E334").Value = $d1_333

$d1_350 = New-Object "object[,]" 3,5
$d1_350[0,0] = '15:31:33'
$d1_350[0,1] = '17:04'
$d1_350[0,2] = '215A_EL PATO'
$d1_350[0,3] = 93
$d1_350[0,4] = 'LP1912'
$d1_350[1,0] = '16:18:55'
$d1_350[1,1] = '17:04'
$d1_350[1,2] = '23_HERNANDEZ'
$d1_350[1,3] = 46
$d1_350[1,4] = 'LP1912'
$d1_350[2,0] = '15:57:48'
$d1_350[2,1] = '17:04'
$d1_350[2,2] = '11_ETCHEVERRY'
$d1_350[2,3] = 67
$d1_350[2,4] = 'LP1912'
$ws1.Range("A350:E352").Value = $d1_350

$d1_362 = New-Object "object[,]" 3,5
$d1_362[0,0] = '15:31:33'
$d1_362[0,1] = '17:21'
$d1_362[0,2] = '26_HERNANDEZ'
$d1_362[0,3] = 110
$d1_362[0,4] = 'LP1912'
$d1_362[1,0] = '17:14:55'
$d1_362[1,1] = '17:21'
$d1_362[1,2] = '10_OLMOS'
$d1_362[1,3] = 7
$d1_362[1,4] = 'LP1912'
$d1_362[2,0] = '16:33:08'
$d1_362[2,1] = '17:21'
$d1_362[2,2] = '16_SANTA ANA'
$d1_362[2,3] = 48
$d1_362[2,4] = 'LP1912'
$ws1.Range("A362:E364").Value = $d1_362

$d1_381 = New-Object "object[,]" 2,5
$d1_381[0,0] = '16:43:37'
$d1_381[0,1] = '17:40'
$d1_381[0,2] = '16_SANTA ANA'
$d1_381[0,3] = 57
$d1_381[0,4] = 'LP1912'
$d1_381[1,0] = '15:57:48'
$d1_381[1,1] = '17:40'
$d1_381[1,2] = '215B_EL PATO'
$d1_381[1,3] = 103
$d1_381[1,4] = 'LP1912'
$ws1.Range("A381:E382").Value = $d1_381

$d1_398 = New-Object "object[,]" 2,5
$d1_398[0,0] = '17:59:03'
$d1_398[0,1] = '18:06'
$d1_398[0,2] = '23_HERNANDEZ'
$d1_398[0,3] = 7
$d1_398[0,4] = 'LP1912'
$d1_398[1,0] = '17:59:03'
$d1_398[1,1] = '18:06'
$d1_398[1,2] = '17_ROMERO'
$d1_398[1,3] = 7
$d1_398[1,4] = 'LP1912'
$ws1.Range("A398:E399").Value = $d1_398

$d1_446 = New-Object "object[,]" 2,5
$d1_446[0,0] = '18:17:05'
$d1_446[0,1] = '19:12'
$d1_446[0,2] = '10_OLMOS'
$d1_446[0,3] = 55
$d1_446[0,4] = 'LP1912'
$d1_446[1,0] = '17:59:03'
$d1_446[1,1] = '19:12'
$d1_446[1,2] = '16_P MOR-SANTA ANA'
$d1_446[1,3] = 73
$d1_446[1,4] = 'LP1912'
$ws1.Range("A446:E447").Value = $d1_446

$d1_479 = New-Object "object[,]" 2,5
$d1_479[0,0] = '18:37:25'
$d1_479[0,1] = '20:00'
$d1_479[0,2] = '17_ROMERO'
$d1_479[0,3] = 83
$d1_479[0,4] = 'LP1912'
$d1_479[1,0] = '19:56:21'
$d1_479[1,1] = '20:00'
$d1_479[1,2] = '14_ABASTO'
$d1_479[1,3] = 4
$d1_479[1,4] = 'LP1912'
$ws1.Range("A479:E480").Value = $d1_479

$d1_494 = New-Object "object[,]" 2,5
$d1_494[0,0] = '19:56:21'
$d1_494[0,1] = '20:23'
$d1_494[0,2] = '215A_EL PATO'
$d1_494[0,3] = 27
$d1_494[0,4] = 'LP1912'
$d1_494[1,0] = '18:37:25'
$d1_494[1,1] = '20:23'
$d1_494[1,2] = '11_ETCHEVERRY'
$d1_494[1,3] = 106
$d1_494[1,4] = 'LP1912'
$ws1.Range("A494:E495").Value = $d1_494

$d1_500 = New-Object "object[,]" 30,5
$d1_500[0,0] = '20:31:05'
$d1_500[0,1] = '20:31'
$d1_500[0,2] = '16_SANTA ANA'
$d1_500[0,3] = 0
$d1_500[0,4] = 'LP1912'
$d1_500[1,0] = '20:31:05'
$d1_500[1,1] = '20:31'
$d1_500[1,2] = '15_ABASTO'
$d1_500[1,3] = 0
$d1_500[1,4] = 'LP1912'
$d1_500[2,0] = '19:42:02'
$d1_500[2,1] = '20:31'
$d1_500[2,2] = '225_GOMEZ'
$d1_500[2,3] = 49
$d1_500[2,4] = 'LP1912'
$d1_500[3,0] = '18:37:25'
$d1_500[3,1] = '20:32'
$d1_500[3,2] = '225_GOMEZ'
$d1_500[3,3] = 115
$d1_500[3,4] = 'LP1912'
$d1_500[4,0] = '18:58:44'
$d1_500[4,1] = '20:35'
$d1_500[4,2] = '14_ABASTO'
$d1_500[4,3] = 97
$d1_500[4,4] = 'LP1912'
$d1_500[5,0] = '19:42:02'
$d1_500[5,1] = '20:39'
$d1_500[5,2] = '11_ETCHEVERRY'
$d1_500[5,3] = 57
$d1_500[5,4] = 'LP1912'
$d1_500[6,0] = '18:51:07'
$d1_500[6,1] = '20:46'
$d1_500[6,2] = '14X44_ABASTO'
$d1_500[6,3] = 115
$d1_500[6,4] = 'LP1912'
$d1_500[7,0] = '18:58:44'
$d1_500[7,1] = '20:48'
$d1_500[7,2] = '14X44_ABASTO'
$d1_500[7,3] = 110
$d1_500[7,4] = 'LP1912'
$d1_500[8,0] = '19:56:21'
$d1_500[8,1] = '20:52'
$d1_500[8,2] = '23_HERNANDEZ'
$d1_500[8,3] = 56
$d1_500[8,4] = 'LP1912'
$d1_500[9,0] = '19:42:02'
$d1_500[9,1] = '20:52'
$d1_500[9,2] = '15_ABASTO'
$d1_500[9,3] = 70
$d1_500[9,4] = 'LP1912'
$d1_500[10,0] = '19:42:02'
$d1_500[10,1] = '20:53'
$d1_500[10,2] = '23_HERNANDEZ'
$d1_500[10,3] = 71
$d1_500[10,4] = 'LP1912'
$d1_500[11,0] = '18:58:44'
$d1_500[11,1] = '20:56'
$d1_500[11,2] = '10_OLMOS'
$d1_500[11,3] = 118
$d1_500[11,4] = 'LP1912'
$d1_500[12,0] = '19:42:02'
$d1_500[12,1] = '20:57'
$d1_500[12,2] = '27_EL RETIRO'
$d1_500[12,3] = 75
$d1_500[12,4] = 'LP1912'
$d1_500[13,0] = '19:56:21'
$d1_500[13,1] = '21:00'
$d1_500[13,2] = '215B_EL PATO'
$d1_500[13,3] = 64
$d1_500[13,4] = 'LP1912'
$d1_500[14,0] = '19:42:02'
$d1_500[14,1] = '21:01'
$d1_500[14,2] = '215B_EL PATO'
$d1_500[14,3] = 79
$d1_500[14,4] = 'LP1912'
$d1_500[15,0] = '19:42:02'
$d1_500[15,1] = '21:04'
$d1_500[15,2] = '84_COLONIA URQUIZA-ESC 49'
$d1_500[15,3] = 82
$d1_500[15,4] = 'LP1912'
$d1_500[16,0] = '20:15:00'
$d1_500[16,1] = '21:16'
$d1_500[16,2] = '84_COLONIA URQUIZA-ESC 49'
$d1_500[16,3] = 61
$d1_500[16,4] = 'LP1912'
$d1_500[17,0] = '19:42:02'
$d1_500[17,1] = '21:21'
$d1_500[17,2] = '26_HERNANDEZ'
$d1_500[17,3] = 99
$d1_500[17,4] = 'LP1912'
$d1_500[18,0] = '19:42:02'
$d1_500[18,1] = '21:23'
$d1_500[18,2] = '10_OLMOS'
$d1_500[18,3] = 101
$d1_500[18,4] = 'LP1912'
$d1_500[19,0] = '20:31:05'
$d1_500[19,1] = '21:31'
$d1_500[19,2] = '84_COLONIA URQUIZA-ESC 49'
$d1_500[19,3] = 60
$d1_500[19,4] = 'LP1912'
$d1_500[20,0] = '20:31:05'
$d1_500[20,1] = '21:33'
$d1_500[20,2] = '23_HERNANDEZ'
$d1_500[20,3] = 62
$d1_500[20,4] = 'LP1912'
$d1_500[21,0] = '19:42:02'
$d1_500[21,1] = '21:38'
$d1_500[21,2] = '14_ABASTO'
$d1_500[21,3] = 116
$d1_500[21,4] = 'LP1912'
$d1_500[22,0] = '19:42:02'
$d1_500[22,1] = '21:38'
$d1_500[22,2] = '17_ROMERO'
$d1_500[22,3] = 116
$d1_500[22,4] = 'LP1912'
$d1_500[23,0] = '20:15:00'
$d1_500[23,1] = '21:43'
$d1_500[23,2] = '17_ROMERO'
$d1_500[23,3] = 88
$d1_500[23,4] = 'LP1912'
$d1_500[24,0] = '19:56:21'
$d1_500[24,1] = '21:47'
$d1_500[24,2] = '215A_EL PATO'
$d1_500[24,3] = 111
$d1_500[24,4] = 'LP1912'
$d1_500[25,0] = '20:31:05'
$d1_500[25,1] = '21:58'
$d1_500[25,2] = '17_ROMERO'
$d1_500[25,3] = 87
$d1_500[25,4] = 'LP1912'
$d1_500[26,0] = '20:15:00'
$d1_500[26,1] = '22:08'
$d1_500[26,2] = '17_ROMERO'
$d1_500[26,3] = 113
$d1_500[26,4] = 'LP1912'
$d1_500[27,0] = '20:31:05'
$d1_500[27,1] = '22:08'
$d1_500[27,2] = '11_ETCHEVERRY'
$d1_500[27,3] = 97
$d1_500[27,4] = 'LP1912'
$d1_500[28,0] = '20:31:05'
$d1_500[28,1] = '22:23'
$d1_500[28,2] = '26_HERNANDEZ'
$d1_500[28,3] = 112
$d1_500[28,4] = 'LP1912'
$d1_500[29,0] = '20:31:05'
$d1_500[29,1] = '22:28'
$d1_500[29,2] = '84_COLONIA URQUIZA-ESC 49'
$d1_500[29,3] = 117
$d1_500[29,4] = 'LP1912'
$ws1.Range("A500:E529").Value = $d1_500

# ---- Sheet LP1912-215 (sheet2) ----
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = 'Última actualización: 20:31:05'

# ---- Sheet 6203-6173 (sheet3) ----
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = 'Última actualización: 20:31:05'
$ws3.Range("A3").Value = 'Total filas: 69'
$d3_74 = New-Object "object[,]" 1,5
$d3_74[0,0] = '20:31:05'
$d3_74[0,1] = '22:21'
$d3_74[0,2] = '215B_LP-P MOR-40 Y 115'
$d3_74[0,3] = 110
$d3_74[0,4] = 'L6173'
$ws3.Range("A74:E74").Value = $d3_74
